$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.564.99"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "'2.556.70"
$ws.Range("E3").Value = "  +5.49%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'573.35"
$ws.Range("E5").Value = "  +2.91%  "
$ws.Range("D6").Value = "'150.57"
$ws.Range("E6").Value = "  +8.63%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "'2.554.44"
$ws.Range("E9").Value = "  +5.45%  "
$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  +2.24%  "
$ws.Range("D11").Value = "'5.75"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "'0.152"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").Value = "'0.358"
$ws.Range("E13").Value = "  +3.36%  "
$ws.Range("D14").Value = "'28.11"
$ws.Range("E14").Value = "  +9.17%  "
$ws.Range("D15").Value = "'3.014.83"
$ws.Range("E15").Value = "  +5.62%  "
$ws.Range("D16").Value = "'63.504.49"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").Value = "'2.539.01"
$ws.Range("E18").Value = "  +4.95%  "
$ws.Range("D19").Value = "'11.61"
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("D20").Value = "'341.80"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "'4.36"
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("D22").Value = "'6.90"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'66.20"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'8.39"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").Value = "'1.47"
$ws.Range("E29").Value = "  +6.93%  "
$ws.Range("D30").Value = "'7.10"
$ws.Range("E30").Value = "  +12.21%  "
$ws.Range("D31").Value = "'0.0₃0839"
$ws.Range("E31").Value = "  +6.37%  "
$ws.Range("D32").Value = "'1.87"
$ws.Range("E32").Value = "  +3.33%  "
$ws.Range("E33").Value = "  +3.91%  "
$ws.Range("D34").Value = "'1.57"
$ws.Range("E34").Value = "  +9.23%  "
$ws.Range("D35").Value = "'413.86"
$ws.Range("E35").Value = "  +10.24%  "
$ws.Range("D36").Value = "'0.406"
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("D37").Value = "'19.10"
$ws.Range("E37").Value = "  +3.12%  "
$ws.Range("D38").Value = "'4.45"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +4.67%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'40.35"
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("D43").Value = "'155.72"
$ws.Range("E43").Value = "  +6.76%  "
$ws.Range("D44").Value = "'3.79"
$ws.Range("E44").Value = "  +3.49%  "
$ws.Range("D45").Value = "'21.05"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").Value = "'0.611"
$ws.Range("E46").Value = "  +4.04%  "
$ws.Range("D47").Value = "'0.0531"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("D48").Value = "'0.0967"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").Value = "'0.0233"
$ws.Range("E49").Value = "  +5.72%  "
$ws.Range("D50").Value = "'18.81"
$ws.Range("E50").Value = "  +4.30%  "
$ws.Range("D51").Value = "'1.87"
$ws.Range("E51").Value = "  +8.82%  "
